# feat: add 2022-Q1 data
#
# Before: sheets = "2021-Q2", "2021-Q4", "总计"
# After:  sheets = "2021-Q2", "2021-Q4", "2022-Q1", "总计"
#
# The old "总计" (Total) worksheet is renamed to "2022-Q1" and repopulated
# with the per-fund holdings table for the new quarter (same shape as the
# "2021-Q2"/"2021-Q4" sheets). A brand-new "总计" worksheet is appended
# after it, carrying the same rollup table as before plus a new top row
# for "2022-Q1".

$wb = $excel.ActiveWorkbook

# Template used to copy header / index-column formatting for the
# per-fund table (same layout as the "2021-Q4" sheet).
$template = $wb.Worksheets.Item("2021-Q4")

# ---------------------------------------------------------------------
# 1. Repurpose the existing "总计" sheet as the new "2022-Q1" sheet.
# ---------------------------------------------------------------------
$q1 = $wb.Worksheets.Item("总计")
$q1.Cells.Clear()
$q1.Name = "2022-Q1"

# Copy header-row formatting (bold, centered, bordered) from 2021-Q4.
$template.Range("B1:H1").Copy()
$q1.Range("B1:H1").PasteSpecial(-4122)

# Copy the formatting used on the index column (A) from 2021-Q4
# (single source cell, repeated across the larger destination range).
$template.Range("A2").Copy()
$q1.Range("A2:A7").PasteSpecial(-4122)

$q1.Range("B1").Value = "基金代码"
$q1.Range("C1").Value = "基金名称"
$q1.Range("D1").Value = "基金规模"
$q1.Range("E1").Value = "股票总仓位"
$q1.Range("F1").Value = "仓位占比"
$q1.Range("G1").Value = "持有市值(亿元)"
$q1.Range("H1").Value = "仓位排名"

$q1Rows = @(
    @(0, "000690", "前海开源大海洋战略经济灵活配置混合", "6.91", "93.87", "4.01", "0.2771", 10),
    @(1, "501201", "红土创新科技创新 3 年封闭运作灵活配置混合", "3.99", "96.70", "4.00", "0.1596", 5),
    @(2, "000969", "前海开源大安全核心精选灵活配置混合", "1.39", "91.04", "4.35", "0.0605", 3),
    @(3, "001060", "前海开源高端装备制造灵活配置混合", "0.97", "89.88", "4.32", "0.0419", 3),
    @(4, "168401", "红土创新转型精选灵活配置混合（LOF）", "0.78", "93.82", "4.46", "0.0348", 3),
    @(5, "350002", "天治低碳经济灵活配置混合", "0.76", "65.23", "2.27", "0.0173", 10)
)

# Fund code / size / position columns are stored as *text* (same shape
# as the "2021-Q2"/"2021-Q4" sheets), not numbers — force text storage
# so leading zeros and trailing-zero decimals round-trip exactly.
$q1.Range("B2:B7").NumberFormat = "@"
$q1.Range("D2:G7").NumberFormat = "@"

$r = 2
foreach ($row in $q1Rows) {
    $q1.Cells.Item($r, 1).Value = $row[0]
    $q1.Cells.Item($r, 2).Value = $row[1]
    $q1.Cells.Item($r, 3).Value = $row[2]
    $q1.Cells.Item($r, 4).Value = $row[3]
    $q1.Cells.Item($r, 5).Value = $row[4]
    $q1.Cells.Item($r, 6).Value = $row[5]
    $q1.Cells.Item($r, 7).Value = $row[6]
    $q1.Cells.Item($r, 8).Value = $row[7]
    $r = $r + 1
}

# ---------------------------------------------------------------------
# 2. Add a fresh "总计" sheet after "2022-Q1" with the rollup table,
#    now including the new 2022-Q1 row on top.
# ---------------------------------------------------------------------
$total = $wb.Worksheets.Add($null, $q1)
$total.Name = "总计"

$template.Range("B1:D1").Copy()
$total.Range("B1:D1").PasteSpecial(-4122)

$template.Range("A2").Copy()
$total.Range("A2:A4").PasteSpecial(-4122)

$total.Range("B1").Value = "日期"
$total.Range("C1").Value = "持有数量(只)"
$total.Range("D1").Value = "持有市值(亿元)"

$totalRows = @(
    @(0, "2022-Q1", 6, 0.59),
    @(1, "2021-Q4", 2, 1.15),
    @(2, "2021-Q2", 2, 0.11)
)

$r = 2
foreach ($row in $totalRows) {
    $total.Cells.Item($r, 1).Value = $row[0]
    $total.Cells.Item($r, 2).Value = $row[1]
    $total.Cells.Item($r, 3).Value = $row[2]
    $total.Cells.Item($r, 4).Value = $row[3]
    $r = $r + 1
}

$wb.Worksheets.Item("2021-Q2").Select()
